$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme
$tcs.Colors(1).RGB = 123456
Write-Host "set done"
